# "added invali cred test case"
# Update the invalid-credentials test rows (rows 4 & 5) on the DATA sheet:
#   username: admin12 -> Admin12
#   password: admin123 -> pass
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

$ws.Range("C4").Value = "Admin12"
$ws.Range("D4").Value = "pass"

$ws.Range("C5").Value = "Admin12"
$ws.Range("D5").Value = "pass"

# Bring the DATA sheet to the front and leave the selection where the
# author's session ended up (D5), matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("D5").Select() | Out-Null
